$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H26").Value = 13500
$ws.Range("I26").Value = 2000
$ws.Range("J26").Value = 25000
$ws.Range("K26").Value = 2000
$ws.Range("L26").Value = 25000
$ws.Range("M26").Value = -1656
$ws.Range("N26").Value = -25688

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H107").Value = 1338.875
$ws.Range("I107").Value = 1226.25
$ws.Range("K107").Value = 1226.25
$ws.Range("M107").Value = 693.75

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H137").Value = 1879.9796
$ws.Range("I137").Value = 1377.1471
$ws.Range("K137").Value = 4131.4413
$ws.Range("M137").Value = -1581.4413

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H138").Value = 3487.848
$ws.Range("I138").Value = 1833.9259
$ws.Range("J138").Value = 5838.1577
$ws.Range("K138").Value = 5501.7777
$ws.Range("L138").Value = 17514.4731
$ws.Range("M138").Value = -361.7776999999996
$ws.Range("N138").Value = -27794.4731

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 1297
$ws.Range("I45").Value = 1176.6111
$ws.Range("J45").Value = 1730.4
$ws.Range("K45").Value = 1176.6111
$ws.Range("L45").Value = 1730.4
$ws.Range("M45").Value = -799.6111000000001
$ws.Range("N45").Value = -2484.4

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 6652.0317
$ws.Range("I61").Value = 4957.434
$ws.Range("J61").Value = 15633.4
$ws.Range("K61").Value = 4957.434
$ws.Range("L61").Value = 15633.4
$ws.Range("M61").Value = -4745.434
$ws.Range("N61").Value = -16057.4

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H62").Value = 38249
$ws.Range("J62").Value = 38249
$ws.Range("L62").Value = 38249
$ws.Range("N62").Value = -39497

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H65").Value = 38249
$ws.Range("J65").Value = 38249
$ws.Range("L65").Value = 114747
$ws.Range("N65").Value = -120987

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 4718.7144
$ws.Range("I74").Value = 1910.6562
$ws.Range("J74").Value = 34671.332
$ws.Range("K74").Value = 1910.6562
$ws.Range("L74").Value = 34671.332
$ws.Range("M74").Value = -1036.6562
$ws.Range("N74").Value = -36419.332

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H75").Value = 32862.5
$ws.Range("J75").Value = 32862.5
$ws.Range("L75").Value = 32862.5
$ws.Range("N75").Value = -34610.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H77").Value = 4718.7144
$ws.Range("I77").Value = 1910.6562
$ws.Range("J77").Value = 34671.332
$ws.Range("K77").Value = 9553.280999999999
$ws.Range("L77").Value = 173356.66
$ws.Range("M77").Value = -5185.280999999999
$ws.Range("N77").Value = -182092.66

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H78").Value = 32862.5
$ws.Range("J78").Value = 32862.5
$ws.Range("L78").Value = 98587.5
$ws.Range("N78").Value = -107323.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H118").Value = 35067.168
$ws.Range("J118").Value = 35067.168
$ws.Range("L118").Value = 35067.168
$ws.Range("N118").Value = -38381.168

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H136").Value = 6652.0317
$ws.Range("I136").Value = 4957.434
$ws.Range("J136").Value = 15633.4
$ws.Range("K136").Value = 14872.302
$ws.Range("L136").Value = 46900.2
$ws.Range("M136").Value = -12322.302
$ws.Range("N136").Value = -52000.2

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H31").Value = 29800
$ws.Range("I31").Value = 0
$ws.Range("J31").Value = 29800
$ws.Range("K31").Value = 0
$ws.Range("L31").Value = 29800
$ws.Range("M31").ClearContents()
$ws.Range("N31").Value = -30304

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H55").Value = 61000
$ws.Range("J55").Value = 61000
$ws.Range("L55").Value = 61000
$ws.Range("N55").Value = -61546

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H69").Value = 36295
$ws.Range("J69").Value = 36295
$ws.Range("L69").Value = 36295
$ws.Range("N69").Value = -37917

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H72").Value = 36295
$ws.Range("J72").Value = 36295
$ws.Range("L72").Value = 108885
$ws.Range("N72").Value = -116997

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H76").Value = 38840.57
$ws.Range("J76").Value = 38840.57
$ws.Range("L76").Value = 38840.57
$ws.Range("N76").Value = -39470.57

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H79").Value = 38840.57
$ws.Range("J79").Value = 38840.57
$ws.Range("L79").Value = 38840.57
$ws.Range("N79").Value = -41024.57

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H21").Value = 0
$ws.Range("J21").Value = 0
$ws.Range("L21").Value = 0
$ws.Range("N21").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1801.9836
$ws.Range("I31").Value = 1252.6274
$ws.Range("J31").Value = 4603.7
$ws.Range("K31").Value = 1252.6274
$ws.Range("L31").Value = 4603.7
$ws.Range("M31").Value = -957.6274000000001
$ws.Range("N31").Value = -5193.7

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 1801.9836
$ws.Range("I34").Value = 1252.6274
$ws.Range("J34").Value = 4603.7
$ws.Range("K34").Value = 1252.6274
$ws.Range("L34").Value = 4603.7
$ws.Range("M34").Value = -1050.6274
$ws.Range("N34").Value = -5007.7

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H94").Value = 1256.5
$ws.Range("I94").Value = 1005
$ws.Range("J94").Value = 1340.3334
$ws.Range("K94").Value = 1005
$ws.Range("L94").Value = 1340.3334
$ws.Range("M94").Value = -554
$ws.Range("N94").Value = -2242.3334

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H42").Value = 2980
$ws.Range("I42").Value = 0
$ws.Range("K42").Value = 0
$ws.Range("M42").ClearContents()

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H113").Value = 696.14734
$ws.Range("I113").Value = 695.0597
$ws.Range("J113").Value = 698.75
$ws.Range("K113").Value = 2085.1791
$ws.Range("L113").Value = 2096.25
$ws.Range("M113").Value = 84.82089999999971
$ws.Range("N113").Value = -6436.25

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H114").Value = 771.13336
$ws.Range("J114").Value = 848.2727
$ws.Range("L114").Value = 2544.8181
$ws.Range("N114").Value = -9052.8181

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 1378.3673
$ws.Range("J131").Value = 1250.75
$ws.Range("L131").Value = 3752.25
$ws.Range("N131").Value = -13832.25

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H46").Value = 26333.334

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H75").Value = 0
$ws.Range("J75").Value = 0
$ws.Range("L75").Value = 0
$ws.Range("N75").ClearContents()

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H78").Value = 0
$ws.Range("J78").Value = 0
$ws.Range("L78").Value = 0
$ws.Range("N78").ClearContents()

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H82").Value = 40328
$ws.Range("J82").Value = 40328
$ws.Range("L82").Value = 40328
$ws.Range("N82").Value = -41094

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H85").Value = 40328
$ws.Range("J85").Value = 40328
$ws.Range("L85").Value = 40328
$ws.Range("N85").Value = -42980

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 4341.027
$ws.Range("I132").Value = 1667.2727
$ws.Range("J132").Value = 26399.5
$ws.Range("K132").Value = 5001.8181
$ws.Range("L132").Value = 79198.5
$ws.Range("M132").Value = -2471.8181
$ws.Range("N132").Value = -84258.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H68").Value = 1101
$ws.Range("I68").Value = 1050
$ws.Range("J68").Value = 1203
$ws.Range("K68").Value = 1050
$ws.Range("L68").Value = 1203
$ws.Range("M68").Value = -301
$ws.Range("N68").Value = -2701

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H71").Value = 1101
$ws.Range("I71").Value = 1050
$ws.Range("J71").Value = 1203
$ws.Range("K71").Value = 5250
$ws.Range("L71").Value = 6015
$ws.Range("M71").Value = -1506
$ws.Range("N71").Value = -13503

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H76").Value = 22884.092
$ws.Range("J76").Value = 22884.092
$ws.Range("L76").Value = 22884.092
$ws.Range("N76").Value = -23560.092

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H79").Value = 22884.092
$ws.Range("J79").Value = 22884.092
$ws.Range("L79").Value = 22884.092
$ws.Range("N79").Value = -25224.092

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H136").Value = 2916.2913
$ws.Range("I136").Value = 1489.2712
$ws.Range("J136").Value = 7126
$ws.Range("K136").Value = 4467.813599999999
$ws.Range("L136").Value = 21378
$ws.Range("M136").Value = -1917.813599999999
$ws.Range("N136").Value = -26478

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H42").Value = 70049
$ws.Range("J42").Value = 70049
$ws.Range("L42").Value = 70049
$ws.Range("N42").Value = -70805

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 2315.05
$ws.Range("I107").Value = 528.2143
$ws.Range("J107").Value = 6484.3335
$ws.Range("K107").Value = 1584.6429
$ws.Range("L107").Value = 19453.0005
$ws.Range("M107").Value = 335.3571000000002
$ws.Range("N107").Value = -23293.0005

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value = 3768.6562
$ws.Range("I113").Value = 5312.4
$ws.Range("J113").Value = 1195.75
$ws.Range("K113").Value = 15937.2
$ws.Range("L113").Value = 3587.25
$ws.Range("M113").Value = -13767.2
$ws.Range("N113").Value = -7927.25
